# Auto-generated edit script: applies the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text assignments (Coin names, Links, Volume%% text, and Price strings that
# are not parseable as a single plain number -- e.g. thousand-grouped "33.999.53" or
# the subscript-digit crypto prices -- these safely round-trip as text on their own).
$ws.Range("D2").Value = "33.999.53"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.780.30"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "2.040.98"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "1.787.37"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E15").Value = "  -1.82%  "
$ws.Range("D16").Value = "33.992.77"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").Value = "1.392.24"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("E44").Value = "  +10.42%  "
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "0.0₆0134"
$ws.Range("E48").Value = "  +8.59%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "1.938.77"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  +0.50%  "

# Price strings that DO look like a plain number (e.g. "0.550", "3.50") would otherwise
# be auto-converted by Excel into a numeric value (dropping trailing zeros / introducing
# floating point noise). Force these to remain literal text, matching the source data,
# by stamping a Text number format for the write and then restoring the default style so
# the cell keeps the workbooks original (unstyled) appearance.
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "226.69"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.550"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "32.35"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.288"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0704"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "10.94"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.621"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "4.13"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "67.84"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "243.53"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.66"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.07"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "160.17"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "16.25"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0511"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.65"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.50"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.80"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.653"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.04"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0187"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.37"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.908"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "77.55"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "13.02"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0497"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "108.10"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "5.83"
$cell.Style = "Normal"
